# Apply the "May 9th" data update: shift the ax..gz (C:H) sensor readings
# up by one row (row N's sensor values move to row N-1) and append 10 new
# rows (timestamps 2000..2900) of walkingToRunning samples, extending the
# used range from A1:H21 to A1:H31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 30,8
$data[0,0] = 0
$data[0,1] = "walkingToRunning"
$data[0,2] = -7.757678985595703
$data[0,3] = -12.53094673156738
$data[0,4] = 0.3439034819602966
$data[0,5] = -0.5661214577783016
$data[0,6] = 0.805141581702479
$data[0,7] = 0.8477507522425696
$data[1,0] = 100
$data[1,1] = "walkingToRunning"
$data[1,2] = -5.902198791503906
$data[1,3] = -8.129295349121094
$data[1,4] = -3.702722311019897
$data[1,5] = 0.1438278629607797
$data[1,6] = 0.8614468347166002
$data[1,7] = 1.115631847037484
$data[2,0] = 200
$data[2,1] = "walkingToRunning"
$data[2,2] = 3.664461135864258
$data[2,3] = -17.59209632873535
$data[2,4] = 3.448472261428833
$data[2,5] = -0.4003904950372941
$data[2,6] = -0.7370264149203714
$data[2,7] = 1.407975565526903
$data[3,0] = 300
$data[3,1] = "walkingToRunning"
$data[3,2] = 12.62621116638184
$data[3,3] = -13.94561958312988
$data[3,4] = 4.861473083496094
$data[3,5] = -1.850929208023019
$data[3,6] = -0.8067298377912074
$data[3,7] = 1.924195033987769
$data[4,0] = 400
$data[4,1] = "walkingToRunning"
$data[4,2] = -1.86917495727539
$data[4,3] = -13.44295883178711
$data[4,4] = -5.517942905426025
$data[4,5] = 1.216999459512474
$data[4,6] = 2.244742146472358
$data[4,7] = 1.425124129069217
$data[5,0] = 500
$data[5,1] = "walkingToRunning"
$data[5,2] = -12.46925640106201
$data[5,3] = -12.17165946960449
$data[5,4] = 5.045797348022461
$data[5,5] = 1.951125894625162
$data[5,6] = -0.8992925304727453
$data[5,7] = -2.80149388804878
$data[6,0] = 600
$data[6,1] = "walkingToRunning"
$data[6,2] = -6.01265811920166
$data[6,3] = -7.838207244873047
$data[6,4] = -2.695719242095948
$data[6,5] = 2.551841642438747
$data[6,6] = -0.6748566922453322
$data[6,7] = -4.086812415073823
$data[7,0] = 700
$data[7,1] = "walkingToRunning"
$data[7,2] = 26.73387336730957
$data[7,3] = -35.77811050415039
$data[7,4] = -2.94611930847168
$data[7,5] = 3.262301676052138
$data[7,6] = 7.791552522133367
$data[7,7] = -4.036898785030718
$data[8,0] = 800
$data[8,1] = "walkingToRunning"
$data[8,2] = -9.356051445007324
$data[8,3] = -7.922670841217041
$data[8,4] = 5.028885364532471
$data[8,5] = -0.5540796260243908
$data[8,6] = 0.137176956098088
$data[8,7] = 3.155596782251743
$data[9,0] = 900
$data[9,1] = "walkingToRunning"
$data[9,2] = 36.33248519897461
$data[9,3] = -82.45571899414062
$data[9,4] = -2.755016803741455
$data[9,5] = -5.821180953192928
$data[9,6] = -0.3604773762299962
$data[9,7] = 6.847419001392469
$data[10,0] = 1000
$data[10,1] = "walkingToRunning"
$data[10,2] = 12.14756202697754
$data[10,3] = 1.323548316955566
$data[10,4] = -4.848138332366943
$data[10,5] = -1.270634695426653
$data[10,6] = 9.376413344722517
$data[10,7] = 6.209869222542664
$data[11,0] = 1100
$data[11,1] = "walkingToRunning"
$data[11,2] = -27.11720466613769
$data[11,3] = -6.166501045227051
$data[11,4] = -2.111397266387939
$data[11,5] = 3.261642898480927
$data[11,6] = 5.511875408211909
$data[11,7] = -1.028492839066049
$data[12,0] = 1200
$data[12,1] = "walkingToRunning"
$data[12,2] = -23.08919525146484
$data[12,3] = -69.49562072753906
$data[12,4] = 16.15834617614746
$data[12,5] = 5.973592050296736
$data[12,6] = 4.099933066318941
$data[12,7] = -1.225896279221969
$data[13,0] = 1300
$data[13,1] = "walkingToRunning"
$data[13,2] = -10.68708801269531
$data[13,3] = -11.96834945678711
$data[13,4] = -0.4801369309425354
$data[13,5] = 1.483463850217993
$data[13,6] = 4.285772220375612
$data[13,7] = -1.806800371900048
$data[14,0] = 1400
$data[14,1] = "walkingToRunning"
$data[14,2] = 20.62446975708008
$data[14,3] = -62.33862686157227
$data[14,4] = 16.18105316162109
$data[14,5] = -3.507514472474761
$data[14,6] = 1.878604103609455
$data[14,7] = 3.185330727665702
$data[15,0] = 1500
$data[15,1] = "walkingToRunning"
$data[15,2] = -75.94793701171875
$data[15,3] = -8.823372840881348
$data[15,4] = -14.5358943939209
$data[15,5] = -6.308687647593385
$data[15,6] = 2.062971848187983
$data[15,7] = 5.004165634666524
$data[16,0] = 1600
$data[16,1] = "walkingToRunning"
$data[16,2] = -7.253152847290039
$data[16,3] = -10.62908172607422
$data[16,4] = -0.7572973966598511
$data[16,5] = -2.792897569150028
$data[16,6] = 2.23644521801742
$data[16,7] = 0.9470839869115668
$data[17,0] = 1700
$data[17,1] = "walkingToRunning"
$data[17,2] = -12.12154960632324
$data[17,3] = -18.81453704833984
$data[17,4] = 0.9114209413528442
$data[17,5] = 2.260217314221203
$data[17,6] = -3.150663144809698
$data[17,7] = -5.311289759026373
$data[18,0] = 1800
$data[18,1] = "walkingToRunning"
$data[18,2] = 2.474559783935547
$data[18,3] = -9.108858108520508
$data[18,4] = 18.17554664611816
$data[18,5] = 4.283928630278282
$data[18,6] = -5.230592103348566
$data[18,7] = -8.01730884964933
$data[19,0] = 1900
$data[19,1] = "walkingToRunning"
$data[19,2] = -10.10338973999023
$data[19,3] = -6.828082084655762
$data[19,4] = 5.437564849853516
$data[19,5] = 0.06544841687701197
$data[19,6] = -3.185547882748625
$data[19,7] = -1.701872078413778
$data[20,0] = 2000
$data[20,1] = "walkingToRunning"
$data[20,2] = 9.688706398010254
$data[20,3] = -28.64854431152344
$data[20,4] = 0.4038746356964111
$data[20,5] = -8.064484859250227
$data[20,6] = -7.12666956174002
$data[20,7] = 2.896797829067585
$data[21,0] = 2100
$data[21,1] = "walkingToRunning"
$data[21,2] = -21.26608657836914
$data[21,3] = 8.973570823669434
$data[21,4] = -19.81748580932617
$data[21,5] = -3.302934056704871
$data[21,6] = -2.102355150832409
$data[21,7] = -0.8136555189938615
$data[22,0] = 2200
$data[22,1] = "walkingToRunning"
$data[22,2] = -51.4796142578125
$data[22,3] = -10.81356239318848
$data[22,4] = -10.4013614654541
$data[22,5] = 3.127995854800499
$data[22,6] = -0.3605925717305229
$data[22,7] = -0.3430995842844817
$data[23,0] = 2300
$data[23,1] = "walkingToRunning"
$data[23,2] = 69.01158905029297
$data[23,3] = -76.07331848144531
$data[23,4] = 28.56989669799805
$data[23,5] = 5.373623056510048
$data[23,6] = -6.720149974233094
$data[23,7] = -0.7005343388036227
$data[24,0] = 2400
$data[24,1] = "walkingToRunning"
$data[24,2] = -14.15320301055908
$data[24,3] = 7.329947471618652
$data[24,4] = -4.596967697143555
$data[24,5] = 3.014144754901412
$data[24,6] = -5.344652057923007
$data[24,7] = -1.741643652473547
$data[25,0] = 2500
$data[25,1] = "walkingToRunning"
$data[25,2] = 66.06742858886719
$data[25,3] = -30.20949363708496
$data[25,4] = 25.3086986541748
$data[25,5] = -2.249343609072485
$data[25,6] = -2.945478901420671
$data[25,7] = 1.552763677134959
$data[26,0] = 2600
$data[26,1] = "walkingToRunning"
$data[26,2] = -70.19232940673828
$data[26,3] = -22.23063087463379
$data[26,4] = -31.12885093688965
$data[26,5] = -3.602893884648981
$data[26,6] = -3.38362657900938
$data[26,7] = -0.5896314640635061
$data[27,0] = 2700
$data[27,1] = "walkingToRunning"
$data[27,2] = 37.95425033569336
$data[27,3] = 1.236392974853516
$data[27,4] = -15.80910873413086
$data[27,5] = -4.196339511379769
$data[27,6] = -2.290577345287686
$data[27,7] = 1.337645951005574
$data[28,0] = 2800
$data[28,1] = "walkingToRunning"
$data[28,2] = -24.9067497253418
$data[28,3] = -28.87722587585449
$data[28,4] = -10.39637756347656
$data[28,5] = 0.4132739001328325
$data[28,6] = -1.036843425527086
$data[28,7] = 6.001093726797151
$data[29,0] = 2900
$data[29,1] = "walkingToRunning"
$data[29,2] = -7.391507625579834
$data[29,3] = -34.15201568603516
$data[29,4] = -12.95433330535889
$data[29,5] = 3.482481982290137
$data[29,6] = 3.476468096074395
$data[29,7] = 0.00721320909322376

$ws.Range("A2:H31").Value = $data

